$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.983.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.642.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.30"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.874.61"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.635.25"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.966.28"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0724"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.91%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.473.51"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.87%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.27%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.55"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.71%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.40"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.784.21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.43%  "
